$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new price row for the latest scrape
$ws.Range("A12").Value = 45812.39347450116
$ws.Range("B12").Value = "EVOWHEY PROTEIN"
$ws.Range("C12").Value = "2Kg"
$ws.Range("D12").Value = "34,90€"

# Match the A11 timestamp style used for the rest of the date column
$ws.Range("A12").NumberFormat = $ws.Range("A11").NumberFormat
